$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "NY.GDP.MKTP.KD.ZG:XKX"
$ws.Range("C1").Value = "NY.GDP.PCAP.CD:XKX"
$ws.Range("D1").Value = "SP.POP.TOTL:XKX"

# Row labels (column A)
$ws.Range("A2").Value = "NY.GDP.MKTP.CD:XKX:cor-value"
$ws.Range("A3").Value = "NY.GDP.MKTP.CD:XKX:p-value"
$ws.Range("A4").Value = "NY.GDP.MKTP.KD.ZG:XKX:cor-value"
$ws.Range("A5").Value = "NY.GDP.MKTP.KD.ZG:XKX:p-value"
$ws.Range("A6").Value = "NY.GDP.PCAP.CD:XKX:cor-value"
$ws.Range("A7").Value = "NY.GDP.PCAP.CD:XKX:p-value"

# Data values
$ws.Range("B2").Value = -0.5538988797221419
$ws.Range("C2").Value = 0.9992848592101515
$ws.Range("D2").Value = 0.9759295041221888

$ws.Range("B3").Value = 0.03986514052192049
$ws.Range("C3").Value = [double]"1.928304747338352e-18"
$ws.Range("D3").Value = [double]"2.666213829249628e-09"

$ws.Range("C4").Value = -0.5513289597721474
$ws.Range("D4").Value = -0.5357789008933723

$ws.Range("C5").Value = 0.04099012598575606
$ws.Range("D5").Value = 0.04830030324290036

$ws.Range("D6").Value = 0.9676736754732521

$ws.Range("D7").Value = [double]"1.536574612054596e-08"

# Styling for header + row-label cells: bold font, thin box border, center/top alignment
$styledRange = $ws.Range("B1:D1,A2:A7")
$styledRange.Font.Bold = $true
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160
$styledRange.Borders.LineStyle = 1
$styledRange.Borders.Weight = 2
